$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: Professional summary paragraph - neutralize language (no bold)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "demographic coding errors affecting 50M voters, developed", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 4: Key Projects "Impact" line - neutralize language (no bold)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 2: Work-experience bullet - replace the phrase with a bold "50M"
# run, keeping the rest of the sentence as plain text runs.
# ---------------------------------------------------------------------
$bullet = $d.Content
$bullet.Find.Execute("race coding errors affecting all Black and Asian-American voters, developed") | Out-Null
$scoped = $d.Range($bullet.Start, $bullet.End)
$scoped.Find.Execute("all Black and Asian-American") | Out-Null
$scoped.Text = "50M"
$scoped.Bold = 1
$scoped.Font.Color = 5258796

# ---------------------------------------------------------------------
# Hunk 3: Move the "Analytics Supervisor - GSD&M" role block (heading +
# 4 paragraphs) so it follows the "Data Products Manager" role block
# instead of preceding it.
# ---------------------------------------------------------------------
$gsdHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Analytics Supervisor - GSD&M*") {
        $gsdHeading = $i
        break
    }
}

$dpmHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Data Products Manager*") {
        $dpmHeading = $i
        break
    }
}

# The GSD&M block is the heading paragraph plus the 4 paragraphs that
# follow it (role line + 3 bullet points), ending right before the
# "Data Products Manager" heading.
$blockStart = $d.Paragraphs($gsdHeading)
$blockEnd = $d.Paragraphs($dpmHeading - 1)
$moveRange = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)
$moveRange.Cut()

# Find the end of the Data Products Manager block (the "57%" paragraph,
# i.e. the paragraph right before the next Heading3 after it).
$dpmHeadingAfterCut = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Data Products Manager*") {
        $dpmHeadingAfterCut = $i
        break
    }
}
$p = $dpmHeadingAfterCut + 1
while ($p -le $d.Paragraphs.Count -and $d.Paragraphs($p).Range.ParagraphStyle.NameLocal -ne "Heading 3") {
    $p = $p + 1
}
$dpmBlockEndIndex = $p - 1
$dpmBlockEnd = $d.Paragraphs($dpmBlockEndIndex)

$insertionPoint = $d.Range($dpmBlockEnd.Range.End, $dpmBlockEnd.Range.End)
$insertionPoint.Paste()

# Restore the Heading3 style on the pasted heading paragraph (Word's
# paste operation re-uses the destination paragraph mark's style for
# the final pasted paragraph, so it needs to be corrected explicitly).
$restoredHeadingIndex = $dpmBlockEndIndex + 1
$d.Paragraphs($restoredHeadingIndex).Style = "Heading 3"

Write-Host "Edits applied successfully"
